$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 193 - this shifts the existing rows 193:198 down to 194:199
# and keeps formatting consistent with the row above (as Excel does natively).
$ws.Rows.Item(193).Insert()

# Populate the newly inserted row 193 with the new record's data.
$ws.Range("A193").Value = 10
$ws.Range("B193").Value = "Vega Modelo de Temuco"
$ws.Range("C193").Value = "La Araucanía"
$ws.Range("D193").Value = 45021
$ws.Range("E193").Value = 9
$ws.Range("F193").Value = "Fruta"
$ws.Range("G193").Value = 100104
$ws.Range("H193").Value = "Frutos de pepita"
$ws.Range("I193").Value = 100104001
$ws.Range("J193").Value = "Granada"
$ws.Range("K193").Value = "Wonderfull"
$ws.Range("L193").Value = "Primera"
$ws.Range("M193").Value = 120
$ws.Range("N193").Value = 24000
$ws.Range("O193").Value = 24000
$ws.Range("P193").Value = 24000
$ws.Range("Q193").Value = "$/bandeja 15 kilos granel"
$ws.Range("R193").Value = "Provincia de Limarí"
$ws.Range("S193").Value = 1600
$ws.Range("T193").Value = 15
